$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.247.58'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +7.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.535.86'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +7.88%  '
$ws.Range('E4').Value = '  +0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.633'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +25.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.993'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.583.38'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.17'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +14.36%  '
$ws.Range('E11').Value = '  +6.85%  '
$ws.Range('E12').Value = '  +6.86%  '
$ws.Range('E13').Value = '  +1.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.016.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +9.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.141.75'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E17').Value = '  +5.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.570.94'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +9.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.94%  '
$ws.Range('E22').Value = '  +7.85%  '
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('E24').Value = '  +7.04%  '
$ws.Range('E25').Value = '  +6.24%  '
$ws.Range('E26').Value = '  +8.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.673.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0828'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '156.84'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.21%  '
$ws.Range('E33').Value = '  +7.42%  '
$ws.Range('E34').Value = '  +7.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  +10.37%  '
$ws.Range('E37').Value = '  +9.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.851'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.86%  '
$ws.Range('E39').Value = '  +11.86%  '
$ws.Range('E40').Value = '  +8.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '35.18'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '291.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +14.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.102'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.625'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.43%  '
$ws.Range('E45').Value = '  +7.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.989'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.760'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +21.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.22'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +14.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.85'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0236'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.10%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.005.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +12.85%  '
